$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H40").Value = 3809.0667
$ws.Range("I40").Value = 2914
$ws.Range("J40").Value = 3988.08
$ws.Range("K40").Value = 2914
$ws.Range("L40").Value = 3988.08
$ws.Range("M40").Value = -2739
$ws.Range("N40").Value = -4338.08
$ws.Range("H99").Value = 389.8889
$ws.Range("I99").Value = 376.25
$ws.Range("J99").Value = 499
$ws.Range("K99").Value = 1128.75
$ws.Range("L99").Value = 1497
$ws.Range("M99").Value = 369.25
$ws.Range("N99").Value = -4493
$ws.Range("H100").Value = 2388.0454
$ws.Range("I100").Value = 1122.091
$ws.Range("J100").Value = 3654
$ws.Range("K100").Value = 1122.091
$ws.Range("L100").Value = 3654
$ws.Range("M100").Value = -581.0909999999999
$ws.Range("N100").Value = -4736
$ws.Range("H116").Value = 20199.445
$ws.Range("I116").Value = 24827.857
$ws.Range("J116").Value = 4000
$ws.Range("K116").Value = 24827.857
$ws.Range("L116").Value = 4000
$ws.Range("M116").Value = -21385.857
$ws.Range("N116").Value = -10884
$ws.Range("H129").Value = 923.7143
$ws.Range("J129").Value = 0
$ws.Range("L129").Value = 0
$ws.Range("H132").Value = 1569.8125
$ws.Range("I132").Value = 1283.0488
$ws.Range("K132").Value = 3849.1464
$ws.Range("M132").Value = -1319.1464
$ws.Range("H137").Value = 3707732.8
$ws.Range("I137").Value = 2782.8
$ws.Range("J137").Value = 5560207.5
$ws.Range("K137").Value = 8348.400000000001
$ws.Range("L137").Value = 16680622.5
$ws.Range("M137").Value = -5798.400000000001
$ws.Range("N137").Value = -16685722.5
$ws.Range("H141").Value = 4690.2354
$ws.Range("I141").Value = 4795.875
$ws.Range("K141").Value = 14387.625
$ws.Range("M141").Value = -9207.625
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H37").Value = 54827
$ws.Range("J37").Value = 69886
$ws.Range("L37").Value = 69886
$ws.Range("N37").Value = -70432
$ws.Range("H44").Value = 53173.332
$ws.Range("J44").Value = 50008
$ws.Range("L44").Value = 50008
$ws.Range("N44").Value = -50984
$ws.Range("H122").Value = 4158.147
$ws.Range("I122").Value = 3611.5417
$ws.Range("K122").Value = 10834.6251
$ws.Range("M122").Value = -8384.625100000001
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H86").Value = 3142.8572
$ws.Range("I86").Value = 2893.375
$ws.Range("J86").Value = 3941.2
$ws.Range("K86").Value = 2893.375
$ws.Range("L86").Value = 3941.2
$ws.Range("M86").Value = -1770.375
$ws.Range("N86").Value = -6187.2
$ws.Range("H89").Value = 3142.8572
$ws.Range("I89").Value = 2893.375
$ws.Range("J89").Value = 3941.2
$ws.Range("K89").Value = 14466.875
$ws.Range("L89").Value = 19706
$ws.Range("M89").Value = -8850.875
$ws.Range("N89").Value = -30938
$ws.Range("H134").Value = 2301398
$ws.Range("I134").Value = 2780239.5
$ws.Range("K134").Value = 8340718.5
$ws.Range("M134").Value = -8338183.5
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H22").Value = 1547.5
$ws.Range("I22").Value = 1807.1
$ws.Range("J22").Value = 249.5
$ws.Range("K22").Value = 1807.1
$ws.Range("L22").Value = 249.5
$ws.Range("M22").Value = -1457.1
$ws.Range("N22").Value = -949.5
$ws.Range("H92").Value = 64985
$ws.Range("J92").Value = 64985
$ws.Range("L92").Value = 64985
$ws.Range("N92").Value = -69977
$ws.Range("H134").Value = 2590.7036
$ws.Range("I134").Value = 2523.9583
$ws.Range("K134").Value = 7571.874899999999
$ws.Range("M134").Value = -5036.874899999999
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H18").Value = 2511.6667
$ws.Range("I18").Value = 1041.2
$ws.Range("K18").Value = 3123.6
$ws.Range("M18").Value = -2954.6
$ws.Range("H107").Value = 757.2727
$ws.Range("J107").Value = 761
$ws.Range("L107").Value = 2283
$ws.Range("N107").Value = -6123
$ws.Range("H120").Value = 15213.421
$ws.Range("I120").Value = 6405.5
$ws.Range("K120").Value = 19216.5
$ws.Range("M120").Value = -14378.5
$ws.Range("H129").Value = 1492.375
$ws.Range("I129").Value = 374.6
$ws.Range("K129").Value = 1123.8
$ws.Range("M129").Value = 3876.2
$ws.Range("H131").Value = 1654.6364
$ws.Range("J131").Value = 1785
$ws.Range("L131").Value = 5355
$ws.Range("N131").Value = -15435
$ws.Range("H140").Value = 1858.3077
$ws.Range("I140").Value = 1350.8889
$ws.Range("J140").Value = 3000
$ws.Range("K140").Value = 4052.6667
$ws.Range("L140").Value = 9000
$ws.Range("M140").Value = 1127.3333
$ws.Range("N140").Value = -19360
$ws.Range("H141").Value = 1500
$ws.Range("I141").Value = 1250
$ws.Range("J141").Value = 2000
$ws.Range("K141").Value = 3750
$ws.Range("L141").Value = 6000
$ws.Range("M141").Value = 1430
$ws.Range("N141").Value = -16360
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H126").Value = 3003.4443
$ws.Range("I126").Value = 2992.6667
$ws.Range("J126").Value = 3008.8333
$ws.Range("K126").Value = 8978.000100000001
$ws.Range("L126").Value = 9026.499899999999
$ws.Range("M126").Value = -6508.000100000001
$ws.Range("N126").Value = -13966.4999
$ws.Range("H132").Value = 2310.923
$ws.Range("I132").Value = 2449.111
$ws.Range("K132").Value = 7347.333
$ws.Range("M132").Value = -4817.333
$ws.Range("H138").Value = 46809.668
$ws.Range("J138").Value = 63619.332
$ws.Range("L138").Value = 63619.332
$ws.Range("N138").Value = -73899.33199999999
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H61").Value = 3563.1365
$ws.Range("I61").Value = 1434.6154
$ws.Range("J61").Value = 6637.6665
$ws.Range("K61").Value = 1434.6154
$ws.Range("L61").Value = 6637.6665
$ws.Range("M61").Value = -1232.6154
$ws.Range("N61").Value = -7041.6665
$ws.Range("H113").Value = 3563.1365
$ws.Range("I113").Value = 1434.6154
$ws.Range("J113").Value = 6637.6665
$ws.Range("K113").Value = 1434.6154
$ws.Range("L113").Value = 6637.6665
$ws.Range("M113").Value = 735.3846000000001
$ws.Range("N113").Value = -10977.6665
$ws.Range("H136").Value = 9226.959999999999
$ws.Range("I136").Value = 9110.666999999999
$ws.Range("K136").Value = 27332.001
$ws.Range("M136").Value = -24782.001
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H81").Value = 4723.5454
$ws.Range("I81").Value = 3641.5
$ws.Range("J81").Value = 5341.857
$ws.Range("K81").Value = 7283
$ws.Range("L81").Value = 10683.714
$ws.Range("M81").Value = -6222
$ws.Range("N81").Value = -12805.714
$ws.Range("H84").Value = 4723.5454
$ws.Range("I84").Value = 3641.5
$ws.Range("J84").Value = 5341.857
$ws.Range("K84").Value = 36415
$ws.Range("L84").Value = 53418.57
$ws.Range("M84").Value = -31111
$ws.Range("N84").Value = -64026.57
$ws.Range("H119").Value = 55576.285
$ws.Range("J119").Value = 55576.285
$ws.Range("L119").Value = 55576.285
$ws.Range("N119").Value = -65252.285
$ws.Range("H136").Value = 2627.4736
$ws.Range("I136").Value = 2330.3572
$ws.Range("K136").Value = 6991.071599999999
$ws.Range("M136").Value = -4441.071599999999

# Remove N129 in ALC (cell deleted in diff)
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("N129").ClearContents()

